$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark left over from a previous
#    editing session.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Collapse the three runs (split apart by a grammar-check proof
#    error marker) that make up the ">>>  your stuff after this
#    line >>>" paragraph back into a single run, and drop the
#    proofErr markers in the process.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    ">>>  your stuff after this line >>>", $false, $false, $false,
    $false, $false, $true, 1, $false,
    ">>>  your stuff after this line >>>", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Add a new paragraph after "Ben changing things up!" containing
#    the text "This is a great subject." formatted in the theme's
#    Accent 5 blue color.
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Ben changing things up!", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterBen = $target.End

# Insert a new paragraph break plus the new text right after the
# "Ben changing things up!" paragraph.
$target.InsertAfter([char]13 + "This is a great subject.")

# Locate the Paragraph object that now holds the freshly inserted
# text so we can apply formatting (using Paragraph.Range keeps the
# paragraph-mark / pPr formatting in sync with the run formatting).
$newParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq ($afterBen + 1)) {
        $newParagraph = $p
    }
}

$colorXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t>This is a great subject.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$newParagraph.Range.InsertXML($colorXml) | Out-Null
